$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$nl = "`n"

# Column E values (Instructional Methods)
$E20 = "- Classroom (Classroom Facilitated Training: 20)" + $nl + "- Didactic Questioning (Classroom Facilitated Training: 20)" + $nl + "- Demonstration (Classroom Facilitated Training: 20)" + $nl + "- Practical (Classroom Facilitated Training: 70)"
$E25 = "- Classroom (Classroom Facilitated Training: 25)" + $nl + "- Didactic Questioning (Classroom Facilitated Training: 25)" + $nl + "- Demonstration (Classroom Facilitated Training: 25)" + $nl + "- Practical (Classroom Facilitated Training: 70)"

# Column G values (Modes of Assessment)
$G_prac5x2 = "- Practical Exam (1:20, 5)" + $nl + "- Practical Exam (1:20, 5)"
$G_writ10_prac10 = "- Written Exam (1:20, 10)" + $nl + "- Practical Exam (1:20, 10)"
$G_writ10x2_prac10x2 = "- Written Exam (1:20, 10)" + $nl + "- Written Exam (1:20, 10)" + $nl + "- Practical Exam (1:20, 10)" + $nl + "- Practical Exam (1:20, 10)"
$G_prac10 = "- Practical Exam (1:20, 10)"

# Row 2 (LU1)
$ws.Range("E2").Value = $E20
$ws.Range("F2").Value = 130
$ws.Range("G2").Value = $G_prac5x2
$ws.Range("H2").Value = 10

# Row 3 (LU2)
$ws.Range("E3").Value = $E20
$ws.Range("F3").Value = 130
$ws.Range("G3").Value = $G_writ10_prac10
$ws.Range("H3").Value = 20

# Row 4 (LU3)
$ws.Range("E4").Value = $E25
$ws.Range("F4").Value = 145
$ws.Range("G4").Value = $G_writ10x2_prac10x2
$ws.Range("H4").Value = 40

# Row 5 (LU4)
$ws.Range("E5").Value = $E25
$ws.Range("F5").Value = 145
$ws.Range("G5").Value = $G_writ10_prac10
$ws.Range("H5").Value = 20

# Row 6 (LU5)
$ws.Range("E6").Value = $E25
$ws.Range("F6").Value = 145
$ws.Range("G6").Value = $G_writ10_prac10
$ws.Range("H6").Value = 20

# Row 7 (LU6)
$ws.Range("E7").Value = $E25
$ws.Range("F7").Value = 145
$ws.Range("G7").Value = $G_prac10
$ws.Range("H7").Value = 10
